# fall 24 week 6 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E2").Value = 10.9

$ws.Range("D3").Value = 10.13
$ws.Range("E3").Value = 10.82

$ws.Range("C4").Value = 9.869999999999999
$ws.Range("E4").Value = 10.59
$ws.Range("F4").Value = 9.74

$ws.Range("B5").Value = 9.1
$ws.Range("C5").Value = 9.140000000000001
$ws.Range("D5").Value = 9.41
$ws.Range("F5").Value = 10.18

$ws.Range("D6").Value = 10.26
$ws.Range("E6").Value = 9.82
$ws.Range("G6").Value = 10.44
$ws.Range("H6").Value = 10.34

$ws.Range("F7").Value = 9.56
$ws.Range("H7").Value = 9.85

$ws.Range("F8").Value = 9.66
$ws.Range("G8").Value = 10.15
$ws.Range("J8").Value = 11.29

$ws.Range("H10").Value = 8.710000000000001

$wb.Save()
